# Change env url / change user credentials
# ------------------------------------------------------------------
# 1) AddUser sheet: insert a "Name" column (new column D) and populate
#    it with the tester names used by the new/renamed test cases.
# 2) Update a handful of test-case titles / data cells on AddUser to
#    reflect the new negative test cases (missing role / phone / mail).
# 3) Edit User sheet: insert a new row (duplicate of the
#    "mustChangePassword" edit test case) above the existing rows.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "AddUser"
# ------------------------------------------------------------------
$addUser = $wb.Worksheets.Item("AddUser")

# Insert a new column before column D ("Email") to hold the user's Name.
$addUser.Columns.Item(4).Insert()
$addUser.Columns.Item(4).ColumnWidth = $addUser.Columns.Item(3).ColumnWidth()

# Header for the new column.
$addUser.Range("D1").Value = "Name"

# Populate the Name column for every data row.
$addUser.Range("D2").Value = "nayra"
$addUser.Range("D3").Value = "nahla"
$addUser.Range("D4").Value = "nahla"
$addUser.Range("D5").Value = "nero"
$addUser.Range("D6").Value = "mostafa"
$addUser.Range("D7").Value = "nayra"
$addUser.Range("D8").Value = "nahla"
$addUser.Range("D9").Value = "nahla"
$addUser.Range("D10").Value = "nero"
$addUser.Range("D11").Value = "mostafa"
$addUser.Range("D12").Value = "nayra"

# Update the titles for the negative-path test cases (rows 2-5).
$addUser.Range("B2").Value = "Verify add new user without Role"
$addUser.Range("B3").Value = "Verify add new user without phone"
$addUser.Range("B4").Value = "Verify add new user without mail and role"
$addUser.Range("B5").Value = "Verify add new user without mail and role"

# Row 11 ("Must Change Password" case) now uses phone "011".
$addUser.Range("F11").Value = "011"

$addUser.Activate()
$addUser.Range("F11").Select()

# ------------------------------------------------------------------
# Sheet "Edit User"
# ------------------------------------------------------------------
$editUser = $wb.Worksheets.Item("Edit User")

# Insert a new row above row 2, duplicating the
# "Verify edit User Password Settings [mustChangePassword Flag]" case
# that used to live in row 5 (now row 6, after the rows below it shift
# down because of the insert).
$editUser.Rows.Item(2).Insert()
$editUser.Range("A6:J6").Copy($editUser.Range("A2:J2"))

$editUser.Activate()
$editUser.Range("A6:XFD6").Select()
